# pwm2speed.xlsx - update PWM -> speed measurement table.
# Default video PWM is 80 instead of 100; changes in NC_Bluetooth sketches.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("24000 Гц Exact")

# New PWM -> time-per-revolution table. PWM goes in column A, measured
# seconds in column B. Existing rows 2-27 (PWM 65-90) get refreshed
# values; rows 28-37 (PWM 91-100) are brand-new measurements extending
# the table now that the default PWM moved from 100 down to 80.
$data = @(
    @(65, 46.6),
    @(66, 36.04),
    @(67, 29.35),
    @(68, 24.58),
    @(69, 21.58),
    @(70, 18.98),
    @(71, 17.16),
    @(72, 15.65),
    @(73, 14.18),
    @(74, 13.2),
    @(75, 12.23),
    @(76, 11.51),
    @(77, 10.77),
    @(78, 10.29),
    @(79, 9.6300000000000008),
    @(80, 9.18),
    @(81, 8.74),
    @(82, 8.36),
    @(83, 8),
    @(84, 7.74),
    @(85, 7.38),
    @(86, 7.11),
    @(87, 6.87),
    @(88, 6.64),
    @(89, 6.39),
    @(90, 6.22),
    @(91, 6.04),
    @(92, 5.93),
    @(93, 5.73),
    @(94, 5.53),
    @(95, 5.44),
    @(96, 5.25),
    @(97, 5.14),
    @(98, 5.01),
    @(99, 4.88),
    @(100, 4.78)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

# Values are now whole-ish seconds rounded to hundredths rather than
# integers -> switch the displayed number format from "0.00" to "0"
# (matches the new cellXfs entry / chart axis format).
$ws.Range("B2:B37").NumberFormat = "0"

# Selection left wherever the author clicked last.
$ws.Range("O1").Select() | Out-Null

Write-Host "pwm2speed table updated"
